$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for D-column cells whose new value parses as a pure number,
# so Excel keeps them as text (matching the source inlineStr cells) instead of
# auto-converting to a numeric type.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("D2").Value = "65.190.79"
$ws.Range("E2").Value = "  +0.79%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("D3").Value = "3.203.40"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("D5").Value = "575.43"
$ws.Range("E5").Value = "  -0.20%  "

$ws.Range("B6").Value = "Solana"
$ws.Range("D6").Value = "167.62"
$ws.Range("E6").Value = "  -2.32%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("D7").Value = "0.592"
$ws.Range("E7").Value = "  -5.48%  "

$ws.Range("B8").Value = "USDC"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("D9").Value = "0.119"
$ws.Range("E9").Value = "  -3.01%  "

$ws.Range("B10").Value = "Toncoin"
$ws.Range("D10").Value = "6.72"
$ws.Range("E10").Value = "  -0.98%  "

$ws.Range("B11").Value = "Cardano"
$ws.Range("D11").Value = "0.392"
$ws.Range("E11").Value = "  +0.56%  "

$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("D12").Value = "3.757.11"
$ws.Range("E12").Value = "  -0.92%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("D13").Value = "0.128"
$ws.Range("E13").Value = "  -0.56%  "

$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("D14").Value = "65.257.64"
$ws.Range("E14").Value = "  +0.84%  "

$ws.Range("B15").Value = "Avalanche"
$ws.Range("D15").Value = "25.70"
$ws.Range("E15").Value = "  -0.42%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("D16").Value = "3.203.36"
$ws.Range("E16").Value = "  -0.85%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("D17").Value = "0.0000158"
$ws.Range("E17").Value = "  -1.04%  "

$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("D18").Value = "413.11"
$ws.Range("E18").Value = "  -0.80%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("D19").Value = "12.85"
$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("D20").Value = "5.34"
$ws.Range("E20").Value = "  -0.99%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("D21").Value = "7.20"
$ws.Range("E21").Value = "  -0.28%  "

$ws.Range("B22").Value = "Dai"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("B23").Value = "Litecoin"
$ws.Range("D23").Value = "69.55"
$ws.Range("E23").Value = "  -1.29%  "

$ws.Range("B24").Value = "Kaspa"
$ws.Range("D24").Value = "0.203"
$ws.Range("E24").Value = "  -2.21%  "

$ws.Range("B25").Value = "Polygon"
$ws.Range("D25").Value = "0.491"
$ws.Range("E25").Value = "  -0.75%  "

$ws.Range("B26").Value = "PEPE"
$ws.Range("D26").Value = "0.0000105"
$ws.Range("E26").Value = "  -4.95%  "

$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("D27").Value = "8.90"
$ws.Range("E27").Value = "  -0.46%  "

$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("D28").Value = "1.01"
$ws.Range("E28").Value = "  +0.61%  "

$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("D29").Value = "1.84"
$ws.Range("E29").Value = "  -1.39%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("D30").Value = "21.47"
$ws.Range("E30").Value = "  -1.90%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("D31").Value = "4.99"
$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("B32").Value = "Aptos"
$ws.Range("D32").Value = "6.40"
$ws.Range("E32").Value = "  -0.32%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("D33").Value = "1.15"
$ws.Range("E33").Value = "  -0.37%  "

$ws.Range("B34").Value = "Monero"
$ws.Range("D34").Value = "156.83"
$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("D35").Value = "1.37"
$ws.Range("E35").Value = "  -1.56%  "

$ws.Range("B36").Value = "Maker"
$ws.Range("D36").Value = "2.730.72"
$ws.Range("E36").Value = "  -2.03%  "

$ws.Range("B37").Value = "Stacks"
$ws.Range("D37").Value = "1.72"
$ws.Range("E37").Value = "  -0.70%  "

$ws.Range("B38").Value = "EnergySwap"
$ws.Range("D38").Value = "24.30"
$ws.Range("E38").Value = "  -4.17%  "

$ws.Range("B39").Value = "Filecoin"
$ws.Range("D39").Value = "4.15"
$ws.Range("E39").Value = "  -1.35%  "

$ws.Range("B40").Value = "Mantle"
$ws.Range("D40").Value = "0.711"
$ws.Range("E40").Value = "  -1.49%  "

$ws.Range("B41").Value = "Hedera"
$ws.Range("D41").Value = "0.0638"
$ws.Range("E41").Value = "  +1.69%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("D42").Value = "5.62"
$ws.Range("E42").Value = "  -3.01%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("D43").Value = "0.0262"
$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("B44").Value = "Bittensor"
$ws.Range("D44").Value = "295.55"
$ws.Range("E44").Value = "  -2.28%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("D45").Value = "21.62"
$ws.Range("E45").Value = "  -1.81%  "

$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "0.0989"
$ws.Range("E47").Value = "  -1.75%  "

$ws.Range("B48").Value = "dogwifhat"
$ws.Range("D48").Value = "1.95"
$ws.Range("E48").Value = "  -10.27%  "

$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "10.47"
$ws.Range("E49").Value = "  +0.52%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "5.78"
$ws.Range("E50").Value = "  -0.79%  "

$ws.Range("B51").Value = "ONDO"
$ws.Range("D51").Value = "0.903"
$ws.Range("E51").Value = "  -2.87%  "
